# Apply the two textual edits described by the diff:
#  1. Merge " " + "286.984.138-89" into a single run " 286.984.138-89"
#     (no visible text change, just a find/replace over that span so the
#     run boundary collapses).
#  2. Insert ", {{FORMA_PAGAMENTO}}" right after "{{VALOR_FINAL}}" and
#     before the trailing period.

$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# Replace " 286.984.138-89" with itself; Word will re-emit the matched
# range as a single run, merging the previously separate " " and
# "286.984.138-89" runs, while leaving the neighbouring runs untouched.
$d.Content.Find.Execute(
    " 286.984.138-89",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " 286.984.138-89",
    2)

# --- Edit 2 -------------------------------------------------------------
# Insert ", {{FORMA_PAGAMENTO}}" before the final period that follows
# "{{VALOR_FINAL}}".
$d.Content.Find.Execute(
    "{{VALOR_FINAL}}.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{{VALOR_FINAL}}, {{FORMA_PAGAMENTO}}.",
    2)
